# Refresh the NATMI ligand-receptor stats (E2:T17) with the new TPM-based
# values (commit: "update scripts wuth new tpm"). Columns A:D (cluster /
# gene labels) are untouched; only the numeric stat columns E..T change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(3.0, 1.0, 10.218698, 30.656094, 0.01131191978527373, 0.01131191978527373, 3.0, 1.0, 40.81054266666667, 122.431628, 0.9943414173631485, 0.9943414173631485, 417.0306107267813, 3753.275496541032, 0.01124791035238732, 0.01124791035238733),
  @(3.0, 1.0, 10.218698, 30.656094, 0.01131191978527373, 0.01131191978527373, 1.0, 0.3333333333333333, 0.01128333333333333, 0.03385, 0.0002749163555820933, 0.0002749163555820933, 0.1153009757666666, 1.0377087819, [double]"3.10983176200443e-06", [double]"3.10983176200443e-06"),
  @(3.0, 1.0, 10.218698, 30.656094, 0.01131191978527373, 0.01131191978527373, 1.0, 0.3333333333333333, 0.1073956666666667, 0.322187, 0.002616675800765965, 0.002616675800765965, 1.097443884175333, 9.876994957577999, [double]"2.95996267623315e-05", [double]"2.959962676233151e-05"),
  @(3.0, 1.0, 10.218698, 30.656094, 0.01131191978527373, 0.01131191978527373, 2.0, 0.6666666666666666, 0.113565, 0.340695, 0.002766990480503436, 0.002766990480503436, 1.16048643837, 10.44437794533, [double]"3.129997436207088e-05", [double]"3.129997436207088e-05"),
  @(3.0, 1.0, 771.4717003333334, 2314.415101, 0.8540056659644313, 0.8540056659644314, 3.0, 1.0, 40.81054266666667, 122.431628, 0.9943414173631485, 0.9943414173631485, 31484.17874257938, 283357.6086832145, 0.8491732043312321, 0.8491732043312322),
  @(3.0, 1.0, 771.4717003333334, 2314.415101, 0.8540056659644313, 0.8540056659644314, 1.0, 0.3333333333333333, 0.01128333333333333, 0.03385, 0.0002749163555820933, 0.0002749163555820933, 8.704772352094444, 78.34295116884999, 0.0002347801253334, 0.0002347801253334),
  @(3.0, 1.0, 771.4717003333334, 2314.415101, 0.8540056659644313, 0.8540056659644314, 1.0, 0.3333333333333333, 0.1073956666666667, 0.322187, 0.002616675800765965, 0.002616675800765965, 82.85271757176523, 745.674458145887, 0.002234655959846149, 0.00223465595984615),
  @(3.0, 1.0, 771.4717003333334, 2314.415101, 0.8540056659644313, 0.8540056659644314, 2.0, 0.6666666666666666, 0.113565, 0.340695, 0.002766990480503436, 0.002766990480503436, 87.61218364835501, 788.5096528351951, 0.002363025548019578, 0.002363025548019578),
  @(3.0, 1.0, 121.2114333333333, 363.6343, 0.1341789346279459, 0.1341789346279459, 3.0, 1.0, 40.81054266666667, 122.431628, 0.9943414173631485, 0.9943414173631485, 4946.704371737823, 44520.3393456404, 0.133419672038229, 0.133419672038229),
  @(3.0, 1.0, 121.2114333333333, 363.6343, 0.1341789346279459, 0.1341789346279459, 1.0, 0.3333333333333333, 0.01128333333333333, 0.03385, 0.0002749163555820933, 0.0002749163555820933, 1.367669006111111, 12.309021055, [double]"3.688798370380282e-05", [double]"3.688798370380283e-05"),
  @(3.0, 1.0, 121.2114333333333, 363.6343, 0.1341789346279459, 0.1341789346279459, 1.0, 0.3333333333333333, 0.1073956666666667, 0.322187, 0.002616675800765965, 0.002616675800765965, 13.01758269045555, 117.1582442141, 0.0003511027712135043, 0.0003511027712135044),
  @(3.0, 1.0, 121.2114333333333, 363.6343, 0.1341789346279459, 0.1341789346279459, 2.0, 0.6666666666666666, 0.113565, 0.340695, 0.002766990480503436, 0.002766990480503436, 13.7653764265, 123.8883878385, 0.0003712718347996191, 0.0003712718347996191),
  @(3.0, 1.0, 0.4548216666666667, 1.364465, 0.0005034796223489374, 0.0005034796223489374, 3.0, 1.0, 40.81054266666667, 122.431628, 0.9943414173631485, 0.9943414173631485, 18.56151903322445, 167.05367129902, 0.0005006306412999052, 0.0005006306412999052),
  @(3.0, 1.0, 0.4548216666666667, 1.364465, 0.0005034796223489374, 0.0005034796223489374, 1.0, 0.3333333333333333, 0.01128333333333333, 0.03385, 0.0002749163555820933, 0.0002749163555820933, 0.005131904472222222, 0.04618714025, [double]"1.384147828860185e-07", [double]"1.384147828860185e-07"),
  @(3.0, 1.0, 0.4548216666666667, 1.364465, 0.0005034796223489374, 0.0005034796223489374, 1.0, 0.3333333333333333, 0.1073956666666667, 0.322187, 0.002616675800765965, 0.002616675800765965, 0.04884587610611111, 0.439612884955, [double]"1.317442943979251e-06", [double]"1.317442943979251e-06"),
  @(3.0, 1.0, 0.4548216666666667, 1.364465, 0.0005034796223489374, 0.0005034796223489374, 2.0, 0.6666666666666666, 0.113565, 0.340695, 0.002766990480503436, 0.002766990480503436, 0.05165182257500001, 0.4648664031750001, [double]"1.393123322166974e-06", [double]"1.393123322166974e-06")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  for ($j = 0; $j -lt $data[$i].Length; $j++) {
    $ws.Cells.Item(2 + $i, 5 + $j).Value = $data[$i][$j]
  }
}

Write-Host "applied"